$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right 50 -> 40, Wrong -1 -> -2, Max text updated
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "38 / 112"
